# ---------------------------------------------------------------------------
# Updates the cryptos price/volume table (Sheet1) with the latest scrape
# values, matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# Columns: A=rank index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many Price cells hold numeric-looking text (e.g. "331.89", "1.894.78") that
# must stay plain text, exactly as authored by the scraper. Assigning such a
# string straight to .Value lets Excel reinterpret it as a real number (losing
# precision / using scientific notation), so we write it with a leading
# apostrophe to force text, then restore a style copied from an already
# unstyled data cell to drop the "quote prefix" formatting Excel adds, keeping
# the cell styleless just like in the original workbook.
$blankStyle = $ws.Range("B2").Style

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = $blankStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "28.918.50"
$ws.Range("E2").Value = "  +1.74%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.889.26"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
Set-TextValue $ws.Range("D5") "331.89"
$ws.Range("E5").Value = "  -1.64%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4626"
$ws.Range("E7").Value = "  -1.80%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.4121"
$ws.Range("E8").Value = "  +3.60%  "

# Row 9
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07986"
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.9960"
$ws.Range("E11").Value = "  -0.35%  "

# Row 12
Set-TextValue $ws.Range("D12") "21.76"
$ws.Range("E12").Value = "  -1.07%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.894.78"
$ws.Range("E13").Value = "  +1.86%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.917"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.066"
$ws.Range("E15").Value = "  -2.61%  "

# Row 16
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
Set-TextValue $ws.Range("D17") "89.18"
$ws.Range("E17").Value = "  -1.54%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001029"
$ws.Range("E18").Value = "  -1.23%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D19") "0.06556"
$ws.Range("E19").Value = "  -1.47%  "

# Row 20
Set-TextValue $ws.Range("D20") "17.49"
$ws.Range("E20").Value = "  -0.42%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.002"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
Set-TextValue $ws.Range("D22") "28.978.57"
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.392"
$ws.Range("E23").Value = "  -1.60%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.22"
$ws.Range("E24").Value = "  +1.40%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.215"
$ws.Range("E25").Value = "  -2.46%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.120.10"
$ws.Range("E26").Value = "  +1.80%  "

# Row 27
Set-TextValue $ws.Range("D27") "157.70"
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
Set-TextValue $ws.Range("D28") "19.70"
$ws.Range("E28").Value = "  -0.19%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.124"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
Set-TextValue $ws.Range("D30") "5.420"
$ws.Range("E30").Value = "  -1.15%  "

# Row 31
Set-TextValue $ws.Range("D31") "117.92"
$ws.Range("E31").Value = "  -1.55%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.9806"
$ws.Range("E32").Value = "  +1.52%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.09379"
$ws.Range("E33").Value = "  -1.45%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.419"
$ws.Range("E34").Value = "  +2.80%  "

# Row 35
$ws.Range("E35").Value = "  +0.35%  "

# Row 36
Set-TextValue $ws.Range("D36") "5.284"
$ws.Range("E36").Value = "  -1.40%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.06076"
$ws.Range("E37").Value = "  -0.67%  "

# Row 38
$ws.Range("E38").Value = "  -0.67%  "

# Row 39
Set-TextValue $ws.Range("D39") "8.325"
$ws.Range("E39").Value = "  +0.12%  "

# Row 40
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.001"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.5777"
$ws.Range("E42").Value = "  -2.71%  "

# Row 43
Set-TextValue $ws.Range("D43") "10.15"
$ws.Range("E43").Value = "  -1.70%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.1820"
$ws.Range("E44").Value = "  -3.04%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.265"
$ws.Range("E45").Value = "  -0.80%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.301"
$ws.Range("E46").Value = "  +11.32%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.5492"
$ws.Range("E47").Value = "  -1.23%  "

# Row 48
Set-TextValue $ws.Range("D48") "12.00"
$ws.Range("E48").Value = "  -0.84%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.906"
$ws.Range("E49").Value = "  -2.46%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.07016"
$ws.Range("E50").Value = "  -2.76%  "

# Row 51
Set-TextValue $ws.Range("D51") "45.90"
$ws.Range("E51").Value = "  +15.77%  "
